# Update gh-pages output (F column "想去人数" counts) to the values
# captured at commit 456a3b4.
#
# Sheet "展览" (rows keyed by F-cell address -> new value)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 293
$ws1.Range("F3").Value = 1179
$ws1.Range("F4").Value = 16727
$ws1.Range("F9").Value = 371
$ws1.Range("F11").Value = 125
$ws1.Range("F12").Value = 11611
$ws1.Range("F14").Value = 1284
$ws1.Range("F15").Value = 4597
$ws1.Range("F16").Value = 425
$ws1.Range("F19").Value = 886

# Sheet "全部类型" (same events repeated, different row offsets)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 293
$ws4.Range("F4").Value = 1179
$ws4.Range("F5").Value = 16727
$ws4.Range("F10").Value = 371
$ws4.Range("F12").Value = 125
$ws4.Range("F15").Value = 11611
$ws4.Range("F17").Value = 1284
$ws4.Range("F18").Value = 4597
$ws4.Range("F19").Value = 425
$ws4.Range("F22").Value = 886
